$wb = $excel.ActiveWorkbook

# --- Update the conversion text on sheet "Hoja1" ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.7 = 6195.74 pesos`n✅ 6195.74 pesos = 1.69 = 913.34 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $newText

# --- Update the rate values on sheet "tasas" ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 587.5
$ws2.Range("O10").Value = 3640
$ws2.Range("N12").Value = 3670
$ws2.Range("O12").Value = 541.01
